# Refined metadata to be additional tab
#
# 1. Update the "time_taken" (column F) timestamps on the existing "data"
#    sheet to their new re-queried values.
# 2. Add a new "metadata" worksheet (placed after "data") describing the
#    PanelApp query that produced the "data" sheet.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh the per-row query timestamps on the "data" sheet ---------
$timeTaken = @{
    2  = "2021-10-05 14:33:36.151155"
    3  = "2021-10-05 14:33:36.151164"
    4  = "2021-10-05 14:33:36.151167"
    5  = "2021-10-05 14:33:36.151170"
    6  = "2021-10-05 14:33:36.151173"
    7  = "2021-10-05 14:33:36.151176"
    8  = "2021-10-05 14:33:36.151178"
    9  = "2021-10-05 14:33:36.151181"
    10 = "2021-10-05 14:33:36.151184"
    11 = "2021-10-05 14:33:36.151187"
    12 = "2021-10-05 14:33:36.151190"
    13 = "2021-10-05 14:33:36.151192"
    14 = "2021-10-05 14:33:36.151195"
    15 = "2021-10-05 14:33:36.151198"
    16 = "2021-10-05 14:33:36.151201"
    17 = "2021-10-05 14:33:36.151203"
    18 = "2021-10-05 14:33:36.151206"
    19 = "2021-10-05 14:33:36.151209"
    20 = "2021-10-05 14:33:36.151212"
    21 = "2021-10-05 14:33:36.151215"
    22 = "2021-10-05 14:33:36.151218"
    23 = "2021-10-05 14:33:36.151220"
    24 = "2021-10-05 14:33:36.151223"
    25 = "2021-10-05 14:33:36.151226"
    26 = "2021-10-05 14:33:36.151229"
    27 = "2021-10-05 14:33:36.151232"
    28 = "2021-10-05 14:33:36.151234"
    29 = "2021-10-05 14:33:36.151237"
    30 = "2021-10-05 14:33:36.151240"
    31 = "2021-10-05 14:33:36.151243"
    32 = "2021-10-05 14:33:36.151246"
    33 = "2021-10-05 14:33:36.151249"
    34 = "2021-10-05 14:33:36.151252"
    35 = "2021-10-05 14:33:36.151254"
    36 = "2021-10-05 14:33:36.151257"
    37 = "2021-10-05 14:33:36.151260"
    38 = "2021-10-05 14:33:36.151263"
    39 = "2021-10-05 14:33:36.151266"
    40 = "2021-10-05 14:33:36.151270"
}

foreach ($row in $timeTaken.Keys) {
    $dataSheet.Range("F$row").Value = $timeTaken[$row]
}

# --- 2. Add the "metadata" worksheet --------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (B1:G1)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row (A2:G2)
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Congenital hypothyroidism"
$metaSheet.Range("C2").Value = 3471

# "data_version" must be stored as text ("0.31"), not a number, so force a
# text number-format on a scratch cell, write the value there, copy it onto
# D2, then clear the scratch cell again.
$scratch = $metaSheet.Range("Z100")
$scratch.NumberFormat = "@"
$scratch.Value = "0.31"
$scratch.Copy()
$metaSheet.Range("D2").PasteSpecial(-4163)
$scratch.Clear()

$metaSheet.Range("E2").Value = "2021-02-12T09:51:20.298166Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:36.147399"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3471/?format=json"

# Match the bold/centered/bordered header style used by the "data" sheet's
# own header row (B1 there carries that style) for the header row and the
# leading index cell A2.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$metaSheet.Range("A2").PasteSpecial(-4122)

$metaSheet.Range("A1").Select()

# Keep "data" as the active/selected tab, matching the original workbook view.
$dataSheet.Activate()
